$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 'AGBCEAAAQBAJ'
$ws.Range("B6").Value = 'Juego de tronos. Realidades, ficciones, turismos'
$ws.Range("C6").Value = 'La influencia de los rodajes cinematográficos o televisivos en el interés turístico de determinados destinos ha ido tomando fuerza en los últimos años. En este libro analizamos el caso del rodaje de una superproducción internacional como es Juego de tronos en localizaciones de Irlanda del Norte o Girona, entre otras, examinando su repercusión en clave turística y los modos en que se entrelazan los atractivos preexistentes con la ficción creada por la productora HBO a partir de los libros de George R. R. Martin. Desde la especificidad del uso turístico de una de las series más exitosas de los últimos años podemos extrapolar diferentes reflexiones relacionadas con fenómenos como la construcción de cánones culturales i de imaginarios, o la noción de autoría en una contemporaneidad en la cual las expectativas del público ganan terreno.'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2021-09-10'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = 'Pere Parramon, F. Xavier Medina'

# Row 7
$ws.Range("A7").Value = 'EqA6DwAAQBAJ'
$ws.Range("B7").Value = 'Juego de tronos y la filosofía'
$ws.Range("C7").Value = 'La lógica es más afilada que las espadas. Se acerca la casa del dragón. Todo lo que siempre has querido saber sobre Juego de tronos, el maravilloso universo creado por R.R. Martin. ¿Son el honor y la verdad necesarios para conseguir la felicidad, o bien nos impiden llegar a ella? ¿Pueden los huargos y otras criaturas fantásticas revelarnos las verdades sobre nuestra conciencia y nuestra realidad? ¿La profecía nos demuestra que somos meros peones del destino o bien que somos libres de vivir una vida auténtica? Si las series de televisión son ideales para el análisis filosófico, Juego de tronos lo es por partida doble. En Westeros y más allá del Mar Angosto, el mundo de George R.R. Martin está repleto de docenas de personajes complejos en conflicto con ellos mismos y en lucha con otros, dudando de sí mismos, abocados al riesgo moral, al engaño, a la incertidumbre, a la arrogancia y a la agitación social y política. Mientras los Siete Reinos están en guerra, más allá del Muro, los horrores del invierno se acercan. Muy lejos, una joven reina lucha con su destino mientras viaja para recuperar su hogar. Todo esto es sabido, pero esta guía perspicaz se basa en las obras de Maquiavelo, Hobbes, Descartes, San Agustín, Platón, Aristóteles y muchos otros grandes filósofos para analizar los personajes y argumentos clave, mientras explora temas como la guerra, el honor, el conocimiento, la moral, la teoría de género y mucho más de una manera tan amena como sorprendente. La crítica ha dicho... «Aplica las teorías de filósofos como Platón, Aristóteles, Kant o Hume para intentar dar explicación a las motivaciones y los conflictos de los personajes de "Juego de tronos". Y que permite a su vez que los lectores se acerquen de una forma curiosa y divertida a esa área del conocimiento humano que es la filosofía.» Fantasymundo «Puede ser la luz que necesitamos para ver y comprender mejor la historia.» The Perks of Being More than a Reader «Sumamente disfrutable.» El Economista «Analiza algunos de los temas claves de la historia de Martin a la luz de los grandes pensadores: honor guerra, conocimiento, moral, verdad...» El norte de Castilla «Me he dado cuenta de muchas cosas de la historia, del porqué de muchos actos de variospersonajes, el arco de unos, las razones de algunas muertes, etc. De pararte a reflexionar con cada línea.» ABIBLIOPHOBICX'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2018-01-18'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = 'Henry Jacoby, William Irwin'

# Row 10
$ws.Range("A10").Value = 'sbPj0AEACAAJ'
$ws.Range("B10").Value = 'Juego de tronos'
$ws.Range("C10").Value = 'Desconocido'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2018'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = 'George R.R Martin'

# Row 11
$ws.Range("A11").Value = 'vfdIEAAAQBAJ'
$ws.Range("B11").Value = 'Juego de tronos. Realidades, ficciones, turismos'
$ws.Range("C11").Value = 'La influencia de los rodajes cinematográficos o televisivos en el interés turístico de determinados destinos ha ido tomando fuerza en los últimos años. En este libro analizamos el caso del rodaje de una superproducción internacional como es Juego de tronos en localizaciones de Irlanda del Norte o Girona, entre otras, examinando su repercusión en clave turística y los modos en que se entrelazan los atractivos preexistentes con la ficción creada por la productora HBO a partir de los libros de George R. R. Martin. Desde la especificidad del uso turístico de una de las series más exitosas de los últimos años podemos extrapolar diferentes reflexiones relacionadas con fenómenos como la construcción de cánones culturales y de imaginarios, o la noción de autoría en una contemporaneidad en la cual las expectativas del público ganan terreno.'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2021-10-11'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = 'Pere Parramon Rubio, Francesc Xavier Medina Luque'
